# Apply the crypto-price refresh described in the commit diff.
# Column D cells that are purely numeric-looking text (e.g. "320.17") get a
# leading apostrophe so Excel stores them as text (quote-prefix), matching
# the original inlineStr/text representation instead of auto-converting them
# to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.656.77'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.857.75'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  -0.76%  '

$ws.Range("D5").Value = '''320.17'
$ws.Range("E5").Value = '  -0.40%  '

$ws.Range("D6").Value = '''1.017'
$ws.Range("E6").Value = '  -0.99%  '

$ws.Range("E7").Value = '  -0.55%  '

$ws.Range("D8").Value = '''0.3790'
$ws.Range("E8").Value = '  +0.37%  '

$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").Value = '''0.8837'
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("D11").Value = '''21.62'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").Value = '1.878.72'
$ws.Range("E12").Value = '  +1.40%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.490'
$ws.Range("E13").Value = '  -0.59%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''6.735'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("D16").Value = '''86.77'
$ws.Range("E16").Value = '  +4.69%  '

$ws.Range("D17").Value = '''1.022'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").Value = '''0.000009064'
$ws.Range("E18").Value = '  +0.33%  '

$ws.Range("D19").Value = '''1.017'
$ws.Range("E19").Value = '  -1.01%  '

$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("D21").Value = '27.659.86'
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").Value = '''5.283'
$ws.Range("E22").Value = '  +0.47%  '

$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").Value = '2.105.09'
$ws.Range("E24").Value = '  +1.70%  '

$ws.Range("D25").Value = '''2.033'
$ws.Range("E25").Value = '  +6.28%  '

$ws.Range("D26").Value = '''157.06'
$ws.Range("E26").Value = '  -0.49%  '

$ws.Range("D27").Value = '''18.70'
$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''1.998'
$ws.Range("E28").Value = '  +1.27%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''5.361'
$ws.Range("E29").Value = '  +1.72%  '

$ws.Range("D30").Value = '''120.45'
$ws.Range("E30").Value = '  +3.02%  '

$ws.Range("D31").Value = '''0.09053'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  +1.72%  '

$ws.Range("D33").Value = '''0.7690'
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("D34").Value = '''3.028'
$ws.Range("E34").Value = '  +5.36%  '

$ws.Range("D35").Value = '''4.556'
$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("D37").Value = '''1.140'
$ws.Range("E37").Value = '  -0.59%  '

$ws.Range("D38").Value = '''0.01975'
$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("D39").Value = '''0.05296'

$ws.Range("D40").Value = '''2.865'
$ws.Range("E40").Value = '  +2.57%  '

$ws.Range("D41").Value = '''0.5192'
$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("D42").Value = '''6.953'
$ws.Range("E42").Value = '  +2.98%  '

$ws.Range("D43").Value = '''0.1680'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").Value = '''8.685'
$ws.Range("E44").Value = '  +2.36%  '

$ws.Range("D45").Value = '''10.79'
$ws.Range("E45").Value = '  +1.98%  '

$ws.Range("D46").Value = '''110.00'
$ws.Range("E46").Value = '  +1.11%  '

$ws.Range("D47").Value = '''1.712'
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").Value = '''0.4723'
$ws.Range("E48").Value = '  +1.50%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.06500'
$ws.Range("E49").Value = '  +1.73%  '

$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '''1.018'
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("D51").Value = '''1.860'
$ws.Range("E51").Value = '  +0.18%  '
